$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mean Expression")

# Insert a new blank column before the existing "CD68+ Membrane PDL1" column (old column E),
# so that: old D (CD8+ Membrane PDL1 data) stays at D, a new blank column is created at E
# (inheriting D's number format/border), and the old E (CD68+ Membrane PDL1 data) shifts to F.
$ws.Columns.Item(5).Insert()

# The freshly inserted column E now has D's style but is empty - move D's original
# ("CD8+ Membrane PDL1" mean-expression numbers) content into E.
$ws.Range("E2").Value2 = $ws.Range("D2").Value2
for ($r = 3; $r -le 8; $r++) {
    $ws.Range("E$r").Value2 = $ws.Range("D$r").Value2
}

# D becomes the new "CD8+ Membrane PD1" column (second marker for the CD8+ phenotype).
$ws.Range("D2").Value2 = "CD8+ Membrane PD1"
$ws.Range("D3").Value2 = 7.434
$ws.Range("D4").Value2 = 6.93463636363636
$ws.Range("D5").Value2 = 2.85891666666667
$ws.Range("D6").Value2 = 3.91146428571429
$ws.Range("D7").Value2 = 5.64085714285714
$ws.Range("D8").Value2 = 5.65262790697674

# Keep the new column the same width as its neighbours.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
